# Add a new row of data to the "Typography" table (row 4) and a new row
# of data to the "Translation" table (row 4), matching the new
# "rotating empty button" entry added to the TouchGFX texts workbook.

$wb = $excel.ActiveWorkbook

# --- Typography sheet: new Typography table row ---
$wsTypo = $wb.Worksheets.Item("Typography")

$wsTypo.Range("B4").Value = "Typography_00"
$wsTypo.Range("C4").Value = "seguisym.ttf"
$wsTypo.Range("D4").Value = 30
$wsTypo.Range("E4").Value = 4
$wsTypo.Range("F4").Value = "?"

# --- Translation sheet: new Translation table row ---
$wsTrans = $wb.Worksheets.Item("Translation")

$wsTrans.Range("B4").Value = "SingleUseId2"
$wsTrans.Range("C4").Value = "Typography_00"
$wsTrans.Range("D4").Value = "Center"
$wsTrans.Range("E4").Value = "LTR"
$wsTrans.Range("F4").Value = "HEJ ERIK"
